$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M9").Value = -202.25
$ws.Range("I9").Value = 371.25
$ws.Range("K9").Value = 371.25
$ws.Range("J9").Value = 2000.6666
$ws.Range("N9").Value = -2338.6666
$ws.Range("L9").Value = 2000.6666
$ws.Range("H9").Value = 815.63635
$ws.Range("H12").Value = 346.7857
$ws.Range("K12").Value = 351
$ws.Range("M12").Value = -181
$ws.Range("I12").Value = 351
$ws.Range("K15").Value = 5766
$ws.Range("H15").Value = 1922
$ws.Range("I15").Value = 1922
$ws.Range("M15").Value = -5597
$ws.Range("J17").Value = 569.51166
$ws.Range("L17").Value = 1708.53498
$ws.Range("N17").Value = -2044.53498
$ws.Range("H17").Value = 569.51166
$ws.Range("J32").Value = 3313.5
$ws.Range("H32").Value = 3313.5
$ws.Range("N32").Value = -3965.5
$ws.Range("L32").Value = 3313.5
$ws.Range("N40").Value = -8249
$ws.Range("I40").Value = 2999
$ws.Range("J40").Value = 7899
$ws.Range("H40").Value = 5449
$ws.Range("M40").Value = -2824
$ws.Range("K40").Value = 2999
$ws.Range("L40").Value = 7899
$ws.Range("N51").Value = -13968
$ws.Range("H51").Value = 19142.857
$ws.Range("J51").Value = 13000
$ws.Range("K51").Value = 20166.666
$ws.Range("I51").Value = 20166.666
$ws.Range("M51").Value = -19682.666
$ws.Range("L51").Value = 13000
$ws.Range("J64").Value = 54099.5
$ws.Range("H64").Value = 54099.5
$ws.Range("L64").Value = 54099.5
$ws.Range("N64").Value = -54595.5
$ws.Range("H67").Value = 54099.5
$ws.Range("J67").Value = 54099.5
$ws.Range("L67").Value = 54099.5
$ws.Range("N67").Value = -55815.5
$ws.Range("M69").Value = -19600.727
$ws.Range("J69").Value = 12774.5
$ws.Range("K69").Value = 20474.727
$ws.Range("N69").Value = -40071.5
$ws.Range("L69").Value = 38323.5
$ws.Range("H69").Value = 9658.048000000001
$ws.Range("I69").Value = 6824.909
$ws.Range("M72").Value = -57056.181
$ws.Range("H72").Value = 9658.048000000001
$ws.Range("K72").Value = 61424.181
$ws.Range("N72").Value = -123706.5
$ws.Range("L72").Value = 114970.5
$ws.Range("I72").Value = 6824.909
$ws.Range("J72").Value = 12774.5
$ws.Range("K74").Value = 9681
$ws.Range("M74").Value = -8745
$ws.Range("H74").Value = 23389.143
$ws.Range("I74").Value = 9681
$ws.Range("H77").Value = 23389.143
$ws.Range("K77").Value = 48405
$ws.Range("I77").Value = 9681
$ws.Range("M77").Value = -43725
$ws.Range("M80").Value = -1225.5001
$ws.Range("H80").Value = 963.46155
$ws.Range("K80").Value = 2223.5001
$ws.Range("I80").Value = 741.1667
$ws.Range("I83").Value = 741.1667
$ws.Range("H83").Value = 963.46155
$ws.Range("M83").Value = -1678.5003
$ws.Range("K83").Value = 6670.5003
$ws.Range("I98").Value = 1183.75
$ws.Range("K98").Value = 1183.75
$ws.Range("M98").Value = 314.25
$ws.Range("H98").Value = 1232.1316
$ws.Range("N116").Value = -18784.625
$ws.Range("L116").Value = 11900.625
$ws.Range("I116").Value = 8582.6
$ws.Range("H116").Value = 10057.277
$ws.Range("J116").Value = 11900.625
$ws.Range("K116").Value = 8582.6
$ws.Range("M116").Value = -5140.6
$ws.Range("K122").Value = 3551.25
$ws.Range("I122").Value = 1183.75
$ws.Range("H122").Value = 1232.1316
$ws.Range("M122").Value = -1101.25
$ws.Range("M132").Value = -20725.889
$ws.Range("K132").Value = 23255.889
$ws.Range("H132").Value = 7250.0645
$ws.Range("I132").Value = 7751.963
$ws.Range("K137").Value = 4369.7142
$ws.Range("M137").Value = -1819.7142
$ws.Range("I137").Value = 1456.5714
$ws.Range("H137").Value = 1985
$ws.Range("M138").Value = -20798
$ws.Range("K138").Value = 25938
$ws.Range("I138").Value = 8646
$ws.Range("H138").Value = 7540.706
$ws.Range("I141").Value = 1339.4445
$ws.Range("L141").Value = 6999
$ws.Range("N141").Value = -17359
$ws.Range("H141").Value = 1587.8334
$ws.Range("M141").Value = 1161.6665
$ws.Range("K141").Value = 4018.3335
$ws.Range("J141").Value = 2333

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L2").Value = 2285.4285
$ws.Range("K2").Value = 794.8333
$ws.Range("M2").Value = -681.8333
$ws.Range("J2").Value = 2285.4285
$ws.Range("I2").Value = 794.8333
$ws.Range("N2").Value = -2511.4285
$ws.Range("H2").Value = 1344
$ws.Range("H32").Value = 36831.676
$ws.Range("K32").Value = 38996.074
$ws.Range("M32").Value = -38709.074
$ws.Range("I32").Value = 38996.074
$ws.Range("H45").Value = 3550.5
$ws.Range("K45").Value = 2278.6
$ws.Range("I45").Value = 2278.6
$ws.Range("M45").Value = -1901.6
$ws.Range("N116").Value = -6873.4285
$ws.Range("L116").Value = 2285.4285
$ws.Range("I116").Value = 794.8333
$ws.Range("H116").Value = 1344
$ws.Range("J116").Value = 2285.4285
$ws.Range("K116").Value = 794.8333
$ws.Range("M116").Value = 1499.1667
$ws.Range("K122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("H122").Value = 4500
$ws.Range("M122").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("M132").Value = -236468
$ws.Range("N132").ClearContents()
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 238998
$ws.Range("L132").Value = 0
$ws.Range("H132").Value = 79666
$ws.Range("I132").Value = 79666

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 794.8333
$ws.Range("N3").Value = -2513.4285
$ws.Range("L3").Value = 2285.4285
$ws.Range("K3").Value = 794.8333
$ws.Range("M3").Value = -680.8333
$ws.Range("H3").Value = 1344
$ws.Range("J3").Value = 2285.4285
$ws.Range("H6").Value = 38541.7
$ws.Range("J6").Value = 38601.89
$ws.Range("N6").Value = -38827.89
$ws.Range("L6").Value = 38601.89
$ws.Range("N20").Value = -4698.8
$ws.Range("J20").Value = 4204.8
$ws.Range("H20").Value = 3848552.5
$ws.Range("L20").Value = 4204.8
$ws.Range("I86").Value = 1638.25
$ws.Range("L86").Value = 4022.6667
$ws.Range("H86").Value = 3068.9
$ws.Range("K86").Value = 1638.25
$ws.Range("N86").Value = -6268.6667
$ws.Range("M86").Value = -515.25
$ws.Range("J86").Value = 4022.6667
$ws.Range("L89").Value = 20113.3335
$ws.Range("K89").Value = 8191.25
$ws.Range("I89").Value = 1638.25
$ws.Range("H89").Value = 3068.9
$ws.Range("J89").Value = 4022.6667
$ws.Range("M89").Value = -2575.25
$ws.Range("N89").Value = -31345.3335
$ws.Range("K105").Value = 2667.4644
$ws.Range("H105").Value = 2695.6128
$ws.Range("I105").Value = 2667.4644
$ws.Range("M105").Value = -920.4643999999998
$ws.Range("M134").Value = -3441.5142
$ws.Range("H134").Value = 2728.3901
$ws.Range("K134").Value = 5976.5142
$ws.Range("I134").Value = 1992.1714

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K7").Value = 97.75
$ws.Range("M7").Value = 15.25
$ws.Range("L7").Value = 444.57144
$ws.Range("J7").Value = 444.57144
$ws.Range("I7").Value = 97.75
$ws.Range("H7").Value = 259.6
$ws.Range("N7").Value = -670.5714399999999
$ws.Range("H16").Value = 2449.3333
$ws.Range("L31").Value = 2500
$ws.Range("H31").Value = 1352.1818
$ws.Range("N31").Value = -3090
$ws.Range("J31").Value = 2500
$ws.Range("J34").Value = 2500
$ws.Range("N34").Value = -2904
$ws.Range("H34").Value = 1352.1818
$ws.Range("L34").Value = 2500
$ws.Range("J50").Value = 43750
$ws.Range("H50").Value = 43750
$ws.Range("N50").Value = -45000
$ws.Range("L50").Value = 43750
$ws.Range("H58").Value = 69220.266
$ws.Range("K58").Value = 85699.25
$ws.Range("I58").Value = 85699.25
$ws.Range("M58").Value = -85496.25
$ws.Range("I62").Value = 5399.5
$ws.Range("H62").Value = 5399.5
$ws.Range("M62").Value = -4775.5
$ws.Range("K62").Value = 5399.5
$ws.Range("M65").Value = -23877.5
$ws.Range("I65").Value = 5399.5
$ws.Range("H65").Value = 5399.5
$ws.Range("K65").Value = 26997.5
$ws.Range("H113").Value = 2449.3333
$ws.Range("M132").Value = -6453.5
$ws.Range("N132").Value = -14053.454
$ws.Range("J132").Value = 2997.818
$ws.Range("K132").Value = 8983.5
$ws.Range("L132").Value = 8993.454000000002
$ws.Range("H132").Value = 2997.3076
$ws.Range("I132").Value = 2994.5
$ws.Range("I136").Value = 85699.25
$ws.Range("H136").Value = 69220.266
$ws.Range("M136").Value = -254547.75
$ws.Range("K136").Value = 257097.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L10").Value = 1671.7896
$ws.Range("N10").Value = -1949.7896
$ws.Range("I10").Value = 246.5
$ws.Range("H10").Value = 503.21738
$ws.Range("K10").Value = 739.5
$ws.Range("J10").Value = 557.2632
$ws.Range("M10").Value = -600.5
$ws.Range("J11").Value = 1462.5
$ws.Range("H11").Value = 690.65
$ws.Range("L11").Value = 4387.5
$ws.Range("N11").Value = -4667.5
$ws.Range("M14").Value = -15469
$ws.Range("H14").Value = 5214
$ws.Range("I14").Value = 5214
$ws.Range("K14").Value = 15642
$ws.Range("H50").Value = 200175
$ws.Range("I50").Value = 196.33333
$ws.Range("M50").Value = -107.99999
$ws.Range("K50").Value = 588.99999
$ws.Range("H53").Value = 200175
$ws.Range("K53").Value = 588.99999
$ws.Range("I53").Value = 196.33333
$ws.Range("M53").Value = -107.99999
$ws.Range("I131").Value = 12343.333
$ws.Range("M131").Value = -31989.999
$ws.Range("H131").Value = 8744.846
$ws.Range("K131").Value = 37029.999
$ws.Range("M132").Value = -10070
$ws.Range("N132").Value = -35059.9997
$ws.Range("J132").Value = 3333.3333
$ws.Range("K132").Value = 12600
$ws.Range("L132").Value = 29999.9997
$ws.Range("H132").Value = 2850
$ws.Range("I132").Value = 1400
$ws.Range("I136").Value = 4185.4
$ws.Range("H136").Value = 4185.4
$ws.Range("M136").Value = -7456.199999999999
$ws.Range("K136").Value = 12556.2
$ws.Range("K140").Value = 6806.117400000001
$ws.Range("H140").Value = 2547.8696
$ws.Range("M140").Value = -1626.117400000001
$ws.Range("I140").Value = 2268.7058

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L102").Value = 5349.625
$ws.Range("N102").Value = -8593.625
$ws.Range("H102").Value = 4599.615
$ws.Range("J102").Value = 5349.625
$ws.Range("I102").Value = 3399.6
$ws.Range("K102").Value = 3399.6
$ws.Range("M102").Value = -1777.6
$ws.Range("K122").Value = 13500
$ws.Range("N122").Value = -24398.5
$ws.Range("J122").Value = 6499.5
$ws.Range("I122").Value = 4500
$ws.Range("H122").Value = 5499.75
$ws.Range("M122").Value = -11050
$ws.Range("L122").Value = 19498.5
$ws.Range("H126").Value = 4227.7085
$ws.Range("L126").Value = 15228.3
$ws.Range("J126").Value = 5076.1
$ws.Range("N126").Value = -20168.3
$ws.Range("M132").Value = -154843.25
$ws.Range("K132").Value = 157373.25
$ws.Range("H132").Value = 48097.91
$ws.Range("I132").Value = 52457.75
$ws.Range("N134").ClearContents()
$ws.Range("J134").Value = 0
$ws.Range("H134").Value = 0
$ws.Range("L134").Value = 0

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 237.66667
$ws.Range("M55").Value = -112.83334
$ws.Range("I55").Value = 285.83334
$ws.Range("K55").Value = 285.83334
$ws.Range("K122").Value = 10039.5879
$ws.Range("I122").Value = 3346.5293
$ws.Range("H122").Value = 3996.1072
$ws.Range("M122").Value = -7589.5879

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K81").Value = 3143.8332
$ws.Range("N81").Value = -12788
$ws.Range("I81").Value = 1571.9166
$ws.Range("M81").Value = -2082.8332
$ws.Range("H81").Value = 2109.2144
$ws.Range("J81").Value = 5333
$ws.Range("L81").Value = 10666
$ws.Range("J84").Value = 5333
$ws.Range("M84").Value = -10415.166
$ws.Range("L84").Value = 53330
$ws.Range("I84").Value = 1571.9166
$ws.Range("H84").Value = 2109.2144
$ws.Range("N84").Value = -63938
$ws.Range("K84").Value = 15719.166
$ws.Range("J106").Value = 99999
$ws.Range("L106").Value = 99999
$ws.Range("H106").Value = 99999
$ws.Range("N106").Value = -102523
$ws.Range("M132").Value = -259005.242
$ws.Range("K132").Value = 261535.242
$ws.Range("H132").Value = 87178.414
$ws.Range("I132").Value = 87178.414
